# Generate Report for Handback
#
# This reproduces the "handback" pass over the localization-status
# workbook:
#   * Status moves from "Ready for handoff" to
#     "Handed back: in sync with en-US" for every language row.
#   * Two new columns get populated per language sheet:
#       F "Latest Target File"   <- same file reference as column A
#       G "Latest Handback File" <- same file reference as column D
#     (with matching hyperlinks, mirroring the existing A/D links)
#   * Column H "Latest Handback DateTime" is stamped with the
#     handback timestamp (per-language, since each language finishes
#     handback at a slightly different moment).

$wb = $excel.ActiveWorkbook

$statusText = "Handed back: in sync with en-US"

# The "Overview" sheet mirrors each language's status in its own
# zh-cn (col B) / de-de (col C) columns - these share the same
# "Ready for handoff" string, so they flip to the new status too.
$overview = $wb.Worksheets.Item("Overview")
foreach ($row in 2..3) {
    $overview.Cells.Item($row, 2).Value2 = $statusText
    $overview.Cells.Item($row, 3).Value2 = $statusText
}

# zh-cn sheet finished handback at 20:46:28, de-de at 20:46:36.
$languageSheets = @(
    @{ Name = "zh-cn"; HandbackTime = "2016-03-21 20:46:28" },
    @{ Name = "de-de"; HandbackTime = "2016-03-21 20:46:36" }
)

foreach ($lang in $languageSheets) {
    $ws = $wb.Worksheets.Item($lang.Name)

    # Worksheet-level Hyperlinks collection is the reliable way to read
    # back an existing link's target address (Range/Cell-level lookup
    # does not surface .Address).
    $linkByCell = @{}
    foreach ($hl in $ws.Hyperlinks) {
        $key = $hl.Range.Row.ToString() + ":" + $hl.Range.Column.ToString()
        $linkByCell[$key] = $hl
    }

    foreach ($row in 2..3) {
        # Status: "Ready for handoff" -> "Handed back: in sync with en-US"
        $ws.Cells.Item($row, 3).Value2 = $statusText

        $sourceLink = $linkByCell[$row.ToString() + ":1"]   # column A - Source File Name
        $handoffLink = $linkByCell[$row.ToString() + ":4"]  # column D - Latest Handoff File

        # column F - Latest Target File (mirrors column A)
        $targetCell = $ws.Cells.Item($row, 6)
        $ws.Hyperlinks.Add($targetCell, $sourceLink.Address, [Type]::Missing, [Type]::Missing, $sourceLink.TextToDisplay) | Out-Null
        $targetCell.Style = "HyperLink"

        # column G - Latest Handback File (mirrors column D)
        $handbackCell = $ws.Cells.Item($row, 7)
        $ws.Hyperlinks.Add($handbackCell, $handoffLink.Address, [Type]::Missing, [Type]::Missing, $handoffLink.TextToDisplay) | Out-Null
        $handbackCell.Style = "HyperLink"

        # column H - Latest Handback DateTime
        $ws.Cells.Item($row, 8).Value2 = $lang.HandbackTime
    }
}
